# Applies the "how partialcells should operate.xlsx" edit:
#  - inserts two new header lines at the top of the sheet (A1, A2)
#  - inserts two new trailing-label columns ("rest of cells" /
#    "cell except gap cap at end") before the final "ZL_U" column in each
#    of the five mini-tables, pushing the existing ZL_U cell two columns right
#  - moves the selection to A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: push the final column of a mini-table two cells to the
#     right, and fill the freed-up pair of cells with the two new labels.
function Insert-RestOfCellsLabels($lastCol, $row) {
    $destCol = [char]([int][char]$lastCol + 2)
    $nextCol = [char]([int][char]$lastCol + 1)

    $srcRange = $ws.Range("$lastCol$row")
    $destRange = $ws.Range("$destCol$row")
    $destRange.Value = $srcRange.Value()

    $ws.Range("$lastCol$row").Value = "rest of cells"
    $ws.Range("$nextCol$row").Value = "cell except gap cap at end"
}

# Table 1 (rows 4-7): last column with ZL_U is G7
Insert-RestOfCellsLabels "G" 7

# Table 2 (rows 10-13): last column with ZL_U is H13
Insert-RestOfCellsLabels "H" 13

# Table 3 (rows 15-18): last column with ZL_U is G18
Insert-RestOfCellsLabels "G" 18

# Table 4 (rows 21-24): last column with ZL_U is H24
Insert-RestOfCellsLabels "H" 24

# Table 5 (rows 26-29): last column with ZL_U is H29
Insert-RestOfCellsLabels "H" 29

# --- new introductory text at the top of the sheet --------------------
$ws.Range("A1").Value = "number of unit cells not under antenna on ONE HALF of substrate: floor((L_sub-L_ant)/2)-1 "
$ws.Range("A2").Value = "(the -1 is because the edge isn't going to behave like a full unit cell - it's a microstrip open rather than a gap cap)"

# --- move the active selection to A3 -----------------------------------
[void]$ws.Range("A3").Select()
